$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---- Update defined name range (D2:E296 -> D2:E307) ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tabelle1!implemented_cards") {
        $n.RefersTo = "=Tabelle1!`$D`$2:`$E`$307"
    }
}

# ---- Bulk-set D2:D383 (existing block), preserving blanks as $null ----
$d2_383 = New-Object 'object[,]' 382,1
$d2_383[0,0] = "Bite"
$d2_383[1,0] = "Claw"
$d2_383[2,0] = "Druid of the Claw"
$d2_383[3,0] = "Healing Touch"
$d2_383[4,0] = "Innervate"
$d2_383[5,0] = "Ironbark Protector"
$d2_383[6,0] = "Keeper of the Grove"
$d2_383[7,0] = "Mark of Nature"
$d2_383[8,0] = "Mark of the Wild"
$d2_383[9,0] = "Moonfire"
$d2_383[10,0] = "Naturalize"
$d2_383[11,0] = "Nourish"
$d2_383[12,0] = "Power of the Wild"
$d2_383[13,0] = "Savage Roar"
$d2_383[14,0] = "Savagery"
$d2_383[15,0] = "Soul of the Forest"
$d2_383[16,0] = "Starfall"
$d2_383[17,0] = "Starfire"
$d2_383[18,0] = "Swipe"
$d2_383[19,0] = "Wild Growth"
$d2_383[20,0] = "Wrath"
$d2_383[21,0] = "Animal Companion"
$d2_383[22,0] = "Arcane Shot"
$d2_383[23,0] = "Deadly Shot"
$d2_383[24,0] = "Eaglehorn Bow"
$d2_383[25,0] = "Explosive Shot"
$d2_383[26,0] = "Explosive Trap"
$d2_383[27,0] = "Flare"
$d2_383[28,0] = "Freezing Trap"
$d2_383[29,0] = "Houndmaster"
$d2_383[30,0] = "Hunter's Mark"
$d2_383[31,0] = "Kill Command"
$d2_383[32,0] = "Misdirection"
$d2_383[33,0] = "Multi-Shot"
$d2_383[34,0] = "Savannah Highmane"
$d2_383[35,0] = "Scavenging Hyena"
$d2_383[36,0] = "Snake Trap"
$d2_383[37,0] = "Snipe"
$d2_383[38,0] = "Starving Buzzard"
$d2_383[39,0] = "Timber Wolf"
$d2_383[40,0] = "Tracking"
$d2_383[41,0] = "Tundra Rhino"
$d2_383[42,0] = "Unleash the Hounds"
$d2_383[43,0] = "Arcane Explosion"
$d2_383[44,0] = "Arcane Intellect"
$d2_383[45,0] = "Arcane Missiles"
$d2_383[46,0] = "Blizzard"
$d2_383[47,0] = "Cone of Cold"
$d2_383[48,0] = "Counterspell"
$d2_383[49,0] = "Ethereal Arcanist"
$d2_383[50,0] = "Fireball"
$d2_383[51,0] = "Flamestrike"
$d2_383[52,0] = "Frostbolt"
$d2_383[53,0] = "Frost Nova"
$d2_383[54,0] = "Ice Barrier"
$d2_383[55,0] = "Ice Lance"
$d2_383[56,0] = "Kirin Tor Mage"
$d2_383[57,0] = "Mana Wyrm"
$d2_383[58,0] = "Mirror Entity"
$d2_383[59,0] = "Mirror Image"
$d2_383[60,0] = "Polymorph"
$d2_383[61,0] = "Pyroblast"
$d2_383[62,0] = "Sorcerer's Apprentice"
$d2_383[63,0] = "Spellbender"
$d2_383[64,0] = "Vaporize"
$d2_383[65,0] = "Water Elemental"
$d2_383[66,0] = "Abomination"
$d2_383[67,0] = "Abusive Sergeant"
$d2_383[68,0] = "Acidic Swamp Ooze"
$d2_383[69,0] = "Acolyte of Pain"
$d2_383[70,0] = "Alarm-o-bot"
$d2_383[71,0] = "Amani Berserker"
$d2_383[72,0] = "Ancient Brewmaster"
$d2_383[73,0] = "Ancient Mage"
$d2_383[74,0] = "Ancient Watcher"
$d2_383[75,0] = "Angry Chicken"
$d2_383[76,0] = "Arcane Golem"
$d2_383[77,0] = "Archmage"
$d2_383[78,0] = "Argent Commander"
$d2_383[79,0] = "Argent Squire"
$d2_383[80,0] = "Azure Drake"
$d2_383[81,0] = "Bloodfen Raptor"
$d2_383[82,0] = "Bloodsail Corsair"
$d2_383[83,0] = "Bloodsail Raider"
$d2_383[84,0] = "Bluegill Warrior"
$d2_383[85,0] = "Booty Bay Bodyguard"
$d2_383[86,0] = "Boulderfist Ogre"
$d2_383[87,0] = "Chillwind Yeti"
$d2_383[88,0] = "Coldlight Oracle"
$d2_383[89,0] = "Coldlight Seer"
$d2_383[90,0] = "Core Hound"
$d2_383[91,0] = "Crazed Alchemist"
$d2_383[92,0] = "Cult Master"
$d2_383[93,0] = "Dalaran Mage"
$d2_383[94,0] = "Dark Iron Dwarf"
$d2_383[95,0] = "Darkscale Healer"
$d2_383[96,0] = "Defender of Argus"
$d2_383[97,0] = "Demolisher"
$d2_383[98,0] = "Dire Wolf Alpha"
$d2_383[99,0] = "Dragonling Mechanic"
$d2_383[100,0] = "Dread Corsair"
$d2_383[101,0] = "Earthen Ring Farseer"
$d2_383[102,0] = "Elven Archer"
$d2_383[103,0] = "Emperor Cobra"
$d2_383[104,0] = "Faerie Dragon"
$d2_383[105,0] = "Fen Creeper"
$d2_383[106,0] = "Flesheating Ghoul"
$d2_383[107,0] = "Frost Elemental"
$d2_383[108,0] = "Frostwolf Grunt"
$d2_383[109,0] = "Frostwolf Warlord"
$d2_383[110,0] = "Gadgetzan Auctioneer"
$d2_383[111,0] = "Gnomish Inventor"
$d2_383[112,0] = "Goldshire Footman"
$d2_383[113,0] = "Grimscale Oracle"
$d2_383[114,0] = "Gurubashi Berserker"
$d2_383[115,0] = "Harvest Golem"
$d2_383[116,0] = "Imp Master"
$d2_383[117,0] = "Injured Blademaster"
$d2_383[118,0] = "Ironbeak Owl"
$d2_383[119,0] = "Ironforge Rifleman"
$d2_383[120,0] = "Ironfur Grizzly"
$d2_383[121,0] = "Jungle Panther"
$d2_383[122,0] = "Knife Juggler"
$d2_383[123,0] = "Kobold Geomancer"
$d2_383[124,0] = "Leper Gnome"
$d2_383[125,0] = "Lightwarden"
$d2_383[126,0] = "Loot Hoarder"
$d2_383[127,0] = "Lord of the Arena"
$d2_383[128,0] = "Mad Bomber"
$d2_383[129,0] = "Magma Rager"
$d2_383[130,0] = "Mana Addict"
$d2_383[131,0] = "Mana Wraith"
$d2_383[132,0] = "Master Swordsmith"
$d2_383[133,0] = "Mind Control Tech"
$d2_383[134,0] = "Mogu'shan Warden"
$d2_383[135,0] = "Murloc Raider"
$d2_383[136,0] = "Murloc Tidecaller"
$d2_383[137,0] = "Murloc Tidehunter"
$d2_383[138,0] = "Nightblade"
$d2_383[139,0] = "Novice Engineer"
$d2_383[140,0] = "Oasis Snapjaw"
$d2_383[141,0] = "Ogre Magi"
$d2_383[142,0] = "Priestess of Elune"
$d2_383[143,0] = "Questing Adventurer"
$d2_383[144,0] = "Raging Worgen"
$d2_383[145,0] = "Raid Leader"
$d2_383[146,0] = "Ravenholdt Assassin"
$d2_383[147,0] = "Razorfen Hunter"
$d2_383[148,0] = "Reckless Rocketeer"
$d2_383[149,0] = "River Crocolisk"
$d2_383[150,0] = "Scarlet Crusader"
$d2_383[151,0] = "Secret Keeper"
$d2_383[152,0] = "Sen'jin Shieldmasta"
$d2_383[153,0] = "Shattered Sun Cleric"
$d2_383[154,0] = "Shieldbearer"
$d2_383[155,0] = "Silverback Patriarch"
$d2_383[156,0] = "Silver Hand Knight"
$d2_383[157,0] = "Silvermoon Guardian"
$d2_383[158,0] = "Southsea Deckhand"
$d2_383[159,0] = "Spellbreaker"
$d2_383[160,0] = "Spiteful Smith"
$d2_383[161,0] = "Stampeding Kodo"
$d2_383[162,0] = "Stonetusk Boar"
$d2_383[163,0] = "Stormpike Commando"
$d2_383[164,0] = "Stormwind Champion"
$d2_383[165,0] = "Stormwind Knight"
$d2_383[166,0] = "Stranglethorn Tiger"
$d2_383[167,0] = "Sunfury Protector"
$d2_383[168,0] = "Sunwalker"
$d2_383[169,0] = "Tauren Warrior"
$d2_383[170,0] = "Thrallmar Farseer"
$d2_383[171,0] = "Venture Co. Mercenary"
$d2_383[172,0] = "Voodoo Doctor"
$d2_383[173,0] = "War Golem"
$d2_383[174,0] = "Windfury Harpy"
$d2_383[175,0] = "Wisp"
$d2_383[176,0] = "Wolfrider"
$d2_383[177,0] = "Worgen Infiltrator"
$d2_383[178,0] = "Young Dragonhawk"
$d2_383[179,0] = "Young Priestess"
$d2_383[180,0] = "Youthful Brewmaster"
$d2_383[181,0] = "Aldor Peacekeeper"
$d2_383[182,0] = "Argent Protector"
$d2_383[183,0] = "Blessed Champion"
$d2_383[184,0] = "Blessing of Kings"
$d2_383[185,0] = "Blessing of Might"
$d2_383[186,0] = "Blessing of Wisdom"
$d2_383[187,0] = "Consecration"
$d2_383[188,0] = "Divine Favor"
$d2_383[189,0] = "Equality"
$d2_383[190,0] = "Eye for an Eye"
$d2_383[191,0] = "Guardian of Kings"
$d2_383[192,0] = "Hammer of Wrath"
$d2_383[193,0] = "Hand of Protection"
$d2_383[194,0] = "Holy Light"
$d2_383[195,0] = "Holy Wrath"
$d2_383[196,0] = "Humility"
$d2_383[197,0] = "Light's Justice"
$d2_383[198,0] = "Noble Sacrifice"
$d2_383[199,0] = "Redemption"
$d2_383[200,0] = "Repentance"
$d2_383[201,0] = "Truesilver Champion"
$d2_383[202,0] = "Auchenai Soulpriest"
$d2_383[203,0] = "Circle of Healing"
$d2_383[204,0] = "Divine Spirit"
$d2_383[205,0] = "Holy Fire"
$d2_383[206,0] = "Holy Nova"
$d2_383[207,0] = "Holy Smite"
$d2_383[208,0] = "Inner Fire"
$d2_383[209,0] = "Lightspawn"
$d2_383[210,0] = "Lightwell"
$d2_383[211,0] = "Mass Dispel"
$d2_383[212,0] = "Mind Blast"
$d2_383[213,0] = "Mind Control"
$d2_383[214,0] = "Mind Vision"
$d2_383[215,0] = "Northshire Cleric"
$d2_383[216,0] = "Power Word: Shield"
$d2_383[217,0] = "Shadow Madness"
$d2_383[218,0] = "Shadow Word: Death"
$d2_383[219,0] = "Shadow Word: Pain"
$d2_383[220,0] = "Silence"
$d2_383[221,0] = "Temple Enforcer"
$d2_383[222,0] = "Thoughtsteal"
$d2_383[223,0] = "Assassinate"
$d2_383[224,0] = "Assassin's Blade"
$d2_383[225,0] = "Backstab"
$d2_383[226,0] = "Betrayal"
$d2_383[227,0] = "Blade Flurry"
$d2_383[228,0] = "Cold Blood"
$d2_383[229,0] = "Conceal"
$d2_383[230,0] = "Deadly Poison"
$d2_383[231,0] = "Defias Ringleader"
$d2_383[232,0] = "Eviscerate"
$d2_383[233,0] = "Fan of Knives"
$d2_383[234,0] = "Headcrack"
$d2_383[235,0] = "Master of Disguise"
$d2_383[236,0] = "Perdition's Blade"
$d2_383[237,0] = "Sap"
$d2_383[238,0] = "Shadowstep"
$d2_383[239,0] = "Shiv"
$d2_383[240,0] = "SI:7 Agent"
$d2_383[241,0] = "Sinister Strike"
$d2_383[242,0] = "Sprint"
$d2_383[243,0] = "Vanish"
$d2_383[244,0] = "Ancestral Healing"
$d2_383[245,0] = "Ancestral Spirit"
$d2_383[246,0] = "Bloodlust"
$d2_383[247,0] = "Dust Devil"
$d2_383[248,0] = "Earth Shock"
$d2_383[249,0] = "Feral Spirit"
$d2_383[250,0] = "Fire Elemental"
$d2_383[251,0] = "Flametongue Totem"
$d2_383[252,0] = "Forked Lightning"
$d2_383[253,0] = "Frost Shock"
$d2_383[254,0] = "Hex"
$d2_383[255,0] = "Lava Burst"
$d2_383[256,0] = "Lightning Bolt"
$d2_383[257,0] = "Lightning Storm"
$d2_383[258,0] = "Mana Tide Totem"
$d2_383[259,0] = "Rockbiter Weapon"
$d2_383[260,0] = "Stormforged Axe"
$d2_383[261,0] = "Totemic Might"
$d2_383[262,0] = "Unbound Elemental"
$d2_383[263,0] = "Windfury"
$d2_383[264,0] = "Windspeaker"
$d2_383[265,0] = "Blood Imp"
$d2_383[266,0] = "Corruption"
$d2_383[267,0] = "Demonfire"
$d2_383[268,0] = "Doomguard"
$d2_383[269,0] = "Drain Life"
$d2_383[270,0] = "Dread Infernal"
$d2_383[271,0] = "Felguard"
$d2_383[272,0] = "Flame Imp"
$d2_383[273,0] = "Hellfire"
$d2_383[274,0] = "Mortal Coil"
$d2_383[275,0] = "Power Overwhelming"
$d2_383[276,0] = "Sacrificial Pact"
$d2_383[277,0] = "Sense Demons"
$d2_383[278,0] = "Shadow Bolt"
$d2_383[279,0] = "Shadowflame"
$d2_383[280,0] = "Siphon Soul"
$d2_383[281,0] = "Soulfire"
$d2_383[282,0] = "Succubus"
$d2_383[283,0] = "Summoning Portal"
$d2_383[284,0] = "Voidwalker"
$d2_383[285,0] = "Arathi Weaponsmith"
$d2_383[286,0] = "Arcanite Reaper"
$d2_383[287,0] = "Armorsmith"
$d2_383[288,0] = "Battle Rage"
$d2_383[289,0] = "Charge"
$d2_383[290,0] = "Cleave"
$d2_383[291,0] = "Commanding Shout"
$d2_383[292,0] = "Cruel Taskmaster"
$d2_383[293,0] = "Execute"
$d2_383[294,0] = "Fiery War Axe"
$d2_383[295,0] = "Frothing Berserker"
$d2_383[296,0] = "Grommash Hellscream"
$d2_383[297,0] = "Heroic Strike"
$d2_383[298,0] = "Inner Rage"
$d2_383[299,0] = "Kor'kron Elite"
$d2_383[300,0] = "Mortal Strike"
$d2_383[301,0] = "Rampage"
$d2_383[302,0] = "Shield Block"
$d2_383[303,0] = "Slam"
$d2_383[304,0] = "Warsong Commander"
$d2_383[305,0] = "Whirlwind"
$d2_383[306,0] = "Core Hound"
$d2_383[307,0] = "Cult Master"
$d2_383[308,0] = "Faerie Dragon"
$d2_383[309,0] = "Gadgetzan Auctioneer"
$d2_383[310,0] = "Stormforged Axe"
$d2_383[311,0] = "PowerOverwhelming"
$d2_383[312,0] = "Soulfire"
$d2_383[313,0] = "Fiery War Axe"
$d2_383[314,0] = "Heroic Strike"
$d2_383[315,0] = "Acolyte of Pain"
$d2_383[316,0] = "Coldlight Oracle"
$d2_383[317,0] = "Flesheating Ghoul"
$d2_383[318,0] = "Mogu'shan Warden"
$d2_383[319,0] = "Eye for an Eye"
$d2_383[320,0] = "Inner Fire"
$d2_383[321,0] = "Thoughtsteal"
$d2_383[322,0] = "Cold Blood"
$d2_383[323,0] = "Conceal"
$d2_383[324,0] = "Blood Imp"
$d2_383[325,0] = "Arathi Weaponsmith"
$d2_383[326,0] = "Northshire Cleric"
$d2_383[327,0] = $null
$d2_383[328,0] = $null
$d2_383[329,0] = "River Crocolisk"
$d2_383[330,0] = "Ice Lance"
$d2_383[331,0] = "Demolisher"
$d2_383[332,0] = "Sen'jin Shieldmasta"
$d2_383[333,0] = $null
$d2_383[334,0] = $null
$d2_383[335,0] = $null
$d2_383[336,0] = $null
$d2_383[337,0] = $null
$d2_383[338,0] = $null
$d2_383[339,0] = $null
$d2_383[340,0] = $null
$d2_383[341,0] = $null
$d2_383[342,0] = $null
$d2_383[343,0] = $null
$d2_383[344,0] = $null
$d2_383[345,0] = $null
$d2_383[346,0] = $null
$d2_383[347,0] = $null
$d2_383[348,0] = $null
$d2_383[349,0] = $null
$d2_383[350,0] = $null
$d2_383[351,0] = $null
$d2_383[352,0] = $null
$d2_383[353,0] = $null
$d2_383[354,0] = $null
$d2_383[355,0] = $null
$d2_383[356,0] = $null
$d2_383[357,0] = $null
$d2_383[358,0] = $null
$d2_383[359,0] = $null
$d2_383[360,0] = $null
$d2_383[361,0] = $null
$d2_383[362,0] = $null
$d2_383[363,0] = $null
$d2_383[364,0] = $null
$d2_383[365,0] = "Counterspell"
$d2_383[366,0] = "Consecration"
$d2_383[367,0] = "Holy Light"
$d2_383[368,0] = "Holy Nova"
$d2_383[369,0] = "Execute"
$d2_383[370,0] = "Gnomish Inventor"
$d2_383[371,0] = "Guardian of Kings"
$d2_383[372,0] = "Abomination"
$d2_383[373,0] = "Divine Spirit"
$d2_383[374,0] = "Druid of the Claw"
$d2_383[375,0] = "Wild Growth"
$d2_383[376,0] = "Freezing Trap"
$d2_383[377,0] = "Ancient Mage"
$d2_383[378,0] = "Knife Juggler"
$d2_383[379,0] = "Circle of Healing"
$d2_383[380,0] = "Lightspawn"
$d2_383[381,0] = "Fan of Knives"
$ws.Range("D2:D383").Value2 = $d2_383

# ---- Clear D356:D366 entirely (value + style) ----
$ws.Range("D356:D366").ClearFormats()
$ws.Range("D356:D366").Value2 = $null

# ---- Apply text style to newly-created D cells (match s="2" format) ----
$ws.Range("D318:D319").NumberFormat = "@"
$ws.Range("D324:D328").NumberFormat = "@"
$ws.Range("D331:D334").NumberFormat = "@"
$ws.Range("D380:D390").NumberFormat = "@"

# ---- Set D384:D390 (brand-new rows) ----
$d384_390 = New-Object 'object[,]' 7,1
$d384_390[0,0] = "Forked Lightning"
$d384_390[1,0] = "Bloodsail Corsair"
$d384_390[2,0] = "Mirror Image"
$d384_390[3,0] = "Ancient Brewmaster"
$d384_390[4,0] = "Emperor Cobra"
$d384_390[5,0] = "Stormpike Commando"
$d384_390[6,0] = "Sap"
$ws.Range("D384:D390").Value2 = $d384_390

